# Update weather metrics (Transavia input automation) for rows 2-6
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("H2").Value = 5.43
$ws.Range("I2").Value = 3.28
$ws.Range("K2").Value = 5.43
$ws.Range("O2").Value = 49
$ws.Range("P2").Value = 0.6899999999999999
$ws.Range("Q2").Value = 71

$ws.Range("H3").Value = 5.44
$ws.Range("I3").Value = 2.89
$ws.Range("J3").Value = 5.44
$ws.Range("O3").Value = 48
$ws.Range("P3").Value = -0.03
$ws.Range("Q3").Value = 72

$ws.Range("H4").Value = 5.55
$ws.Range("I4").Value = 2.99
$ws.Range("J4").Value = 5.55
$ws.Range("L4").Value = 1017
$ws.Range("M4").Value = 1017
$ws.Range("P4").Value = -0.19
$ws.Range("Q4").Value = 75

$ws.Range("H5").Value = 5.89
$ws.Range("I5").Value = 3.04
$ws.Range("J5").Value = 5.89
$ws.Range("L5").Value = 1016
$ws.Range("M5").Value = 1016
$ws.Range("O5").Value = 46
$ws.Range("P5").Value = -0.3
$ws.Range("Q5").Value = 88

$ws.Range("H6").Value = 5.96
$ws.Range("I6").Value = 2.88
$ws.Range("J6").Value = 5.96
$ws.Range("P6").Value = -0.13
$ws.Range("Q6").Value = 94
